$wb = $excel.ActiveWorkbook

# --- Sheet "3a - Uncertainty analysis" (sheet4.xml) ---
# Add a new "Reconstructed motions" boolean column (B), mirroring the same
# column already present on the other sheets ("2a", "4 - Slab pull optimisation", etc).
$ws3a = $wb.Worksheets.Item("3a - Uncertainty analysis")
$ws3a.Range("B1").Value = "Reconstructed motions"
$ws3a.Range("B2").Value = $true
$ws3a.Range("B3").Value = $true
$ws3a.Range("B4").Value = $false
$ws3a.Columns.Item(2).ColumnWidth = 19.330729166666668

# --- Sheet "3b - Uncertainty analysis" (sheet5.xml) ---
$ws3b = $wb.Worksheets.Item("3b - Uncertainty analysis")
$ws3b.Range("B1").Value = "Reconstructed motions"
$ws3b.Range("B2").Value = $true
$ws3b.Range("B3").Value = $false
$ws3b.Columns.Item(2).ColumnWidth = 19.330729166666668

# --- Selection / active-tab bookkeeping ---
# Leave a lingering selection on 3a at B7 (selecting it makes it active
# momentarily; 3b is activated afterwards so it ends up the final active tab).
$ws3a.Range("B7").Select()

# 3b ("3b - Uncertainty analysis") becomes the active sheet/tab, with the
# cursor left on B4 - this also clears tabSelected from whatever sheet was
# previously active (2b - Extended torque balance).
$ws3b.Range("B4").Select()
